$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains its original text formatting so that
# numeric-looking values (e.g. "1.00", "596.58") are not silently converted
# to real numbers by Excel.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D11","D12","D13","D14","D15","D17","D18","D20","D21","D23","D24","D25","D26","D27","D29","D30","D31","D32","D34","D35","D36","D37","D38","D40","D41","D42","D44","D45","D47","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "72.967.02"
$ws.Range("E2").Value = "  +2.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.978.02"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "596.58"
$ws.Range("E5").Value = "  +10.95%  "

# Row 6 - Solana
$ws.Range("D6").Value = "159.69"

# Row 7 - XRP
$ws.Range("D7").Value = "0.681"
$ws.Range("E7").Value = "  -0.63%  "

# Row 8 - USDC
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.13%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.748"
$ws.Range("E9").Value = "  +1.52%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.74%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "53.43"
$ws.Range("E11").Value = "  -3.37%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "0.0000318"
$ws.Range("E12").Value = "  +0.58%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "10.95"
$ws.Range("E13").Value = "  +3.06%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.615.99"
$ws.Range("E14").Value = "  +0.80%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.987.44"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +8.03%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "14.04"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "20.29"
$ws.Range("E18").Value = "  -1.44%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.28%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "72.611.87"
$ws.Range("E20").Value = "  +2.35%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "434.49"
$ws.Range("E21").Value = "  +2.13%  "

# Row 22 - PancakeSwap
$ws.Range("E22").Value = "  +13.00%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "95.84"
$ws.Range("E23").Value = "  -1.17%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "3.41"
$ws.Range("E24").Value = "  -4.41%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "14.17"
$ws.Range("E25").Value = "  -1.55%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "4.34"
$ws.Range("E26").Value = "  +15.19%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "11.19"
$ws.Range("E27").Value = "  -1.59%  "

# Row 28 - LEO
$ws.Range("E28").Value = "  +0.88%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "10.43"
$ws.Range("E29").Value = "  -1.82%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "36.23"
$ws.Range("E30").Value = "  -0.38%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.80"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "13.71"
$ws.Range("E32").Value = "  +2.84%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.21%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "47.78"
$ws.Range("E34").Value = "  -4.65%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "664.62"
$ws.Range("E35").Value = "  -2.71%  "

# Row 36 - OKB
$ws.Range("D36").Value = "70.73"
$ws.Range("E36").Value = "  +8.78%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "0.0₃0898"
$ws.Range("E37").Value = "  +9.53%  "

# Row 38 - TheGraph
$ws.Range("D38").Value = "0.436"
$ws.Range("E38").Value = "  -0.16%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  -0.08%  "

# Row 40 - was Kaspa, now ThetaToken
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "3.34"
$ws.Range("E40").Value = "  -1.48%  "

# Row 41 - was ThetaToken, now WEMIXToken
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "3.33"
$ws.Range("E41").Value = "  +4.54%  "

# Row 42 - was WEMIXToken, now Kaspa
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.145"
$ws.Range("E42").Value = "  -3.21%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.26%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0490"
$ws.Range("E44").Value = "  +2.00%  "

# Row 45 - THORChain
$ws.Range("D45").Value = "10.59"
$ws.Range("E45").Value = "  +7.02%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  +0.34%  "

# Row 47 - ApeXProtocol
$ws.Range("D47").Value = "3.44"
$ws.Range("E47").Value = "  +3.79%  "

# Row 48 - Fetch.AI
$ws.Range("E48").Value = "  -3.38%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.875.78"
$ws.Range("E49").Value = "  +9.49%  "

# Row 50 - Stacks
$ws.Range("E50").Value = "  +0.89%  "

# Row 51 - LidoDAOToken
$ws.Range("E51").Value = "  +4.05%  "
